$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3351
$ws1.Range("F4").Value = 62
$ws1.Range("F5").Value = 1454
$ws1.Range("F6").Value = 32
$ws1.Range("F7").Value = 320

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3351
$ws4.Range("F4").Value = 62
$ws4.Range("F5").Value = 1454
$ws4.Range("F6").Value = 32
$ws4.Range("F8").Value = 320
